$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "'07/18/2025"
$ws.Range("B42").Value = 520.4100000000035
$ws.Range("C42").Value = 0.09607809227339917
$ws.Range("D42").Value = 50
